$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98 (shifts existing rows 98:130 down to 99:131)
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new record
$ws.Range("A98").Value = 10
$ws.Range("B98").Value = "Vega Modelo de Temuco"
$ws.Range("C98").Value = "La Araucanía"
$ws.Range("D98").Value = 44798
$ws.Range("E98").Value = 9
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100107
$ws.Range("H98").Value = "Otros"
$ws.Range("I98").Value = 100107002
$ws.Range("J98").Value = "Chirimoya"
$ws.Range("K98").Value = "Cultivar IV Región"
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 50
$ws.Range("N98").Value = 4500
$ws.Range("O98").Value = 4500
$ws.Range("P98").Value = 4500
$ws.Range("Q98").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R98").Value = "Provincia del Elquí"
$ws.Range("S98").Value = 4500
$ws.Range("T98").Value = 1
